$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data values on row 2 (trial 1)
$ws.Range("E2").Value = 7
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 13

# Update the active selection to match the saved view state
$ws.Range("E2").Select()
